$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (CP040557): add ncbi-spuid_namespace value ---
$ws.Range("C3").Value = "BBL-BIOINTEL"

# --- Row 4: rename sample from BX248355 to BX248355.1-segment2, add namespace ---
$ws.Range("A4").Value = "BX248355.1-segment2"
$ws.Range("B4").Value = "BX248355.1-segment2"
$ws.Range("C4").Value = "BBL-BIOINTEL"
$ws.Range("K4").Value = "BX248355.1-segment2"
$ws.Range("AE4").Value = "BX248355.1-segment2_Cd"

# --- Insert two new rows (5 and 6), inheriting row 4's formatting ---
$ws.Rows("5").Insert()
$ws.Rows("6").Insert()

# --- Row 5: BX248355.1-segment3, cloned from row 4 pattern ---
$ws.Range("A5").Value = "BX248355.1-segment3"
$ws.Range("B5").Value = "BX248355.1-segment3"
$ws.Range("C5").Value = "BBL-BIOINTEL"
$ws.Range("E5").Value = "John Doe"
$ws.Range("F5").Value = "Black Bird Labs"
$ws.Range("G5").Value = "Bio intelligence"
$ws.Range("K5").Value = "BX248355.1-segment3"
$ws.Range("L5").Value = "Clinical"
$ws.Range("O5").Value = "Homo sapiens"
$ws.Range("P5").Value = "Cdiphtheriae"
$ws.Range("Q5").Value = "2022-06"
$ws.Range("R5").Value = "USA"
$ws.Range("T5").Value = "Not provided"
$ws.Range("V5").Value = "Not provided"
$ws.Range("AE5").Value = "BX248355.1-segment3_Cd"
$ws.Range("AL5").Value = "local"
$ws.Range("AM5").Value = "/scicomp/instruments-pure/23-4-631_Illumina-MiSeq-M03083/220622_M03083_0094_000000000-KG73J/Alignment_1/20220624_225908/Fastq/BX248355_R1.fastq.gz"
$ws.Range("AN5").Value = "/scicomp/instruments-pure/23-4-631_Illumina-MiSeq-M03083/220622_M03083_0094_000000000-KG73J/Alignment_1/20220624_225908/Fastq/BX248355_R2.fastq.gz"

# --- Row 6: BX248355.1-segment4, cloned from row 4 pattern ---
$ws.Range("A6").Value = "BX248355.1-segment4"
$ws.Range("B6").Value = "BX248355.1-segment4"
$ws.Range("C6").Value = "BBL-BIOINTEL"
$ws.Range("E6").Value = "John Doe"
$ws.Range("F6").Value = "Black Bird Labs"
$ws.Range("G6").Value = "Bio intelligence"
$ws.Range("K6").Value = "BX248355.1-segment4"
$ws.Range("L6").Value = "Clinical"
$ws.Range("O6").Value = "Homo sapiens"
$ws.Range("P6").Value = "Cdiphtheriae"
$ws.Range("Q6").Value = "2022-06"
$ws.Range("R6").Value = "USA"
$ws.Range("T6").Value = "Not provided"
$ws.Range("V6").Value = "Not provided"
$ws.Range("AE6").Value = "BX248355.1-segment4_Cd"
$ws.Range("AL6").Value = "local"
$ws.Range("AM6").Value = "/scicomp/instruments-pure/23-4-631_Illumina-MiSeq-M03083/220622_M03083_0094_000000000-KG73J/Alignment_1/20220624_225908/Fastq/BX248355_R1.fastq.gz"
$ws.Range("AN6").Value = "/scicomp/instruments-pure/23-4-631_Illumina-MiSeq-M03083/220622_M03083_0094_000000000-KG73J/Alignment_1/20220624_225908/Fastq/BX248355_R2.fastq.gz"

# --- Update the active selection to A4 (also clears any stale topLeftCell scroll state) ---
$ws.Range("A4").Select()
